# Fix the "Pos*5" calculated column in Table1: the multiplier was
# (5 - 5/4) = 3.75 which skipped/offset every entry by the B/C column
# width; it should be (5 - 5/5) = 4 so Pos*5 lines up correctly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo  = $ws.ListObjects.Item("Table1")
$col = $lo.ListColumns.Item("Pos*5")

# Re-write the whole calculated column's formula (this also fixes row 12,
# which previously had a stray non-conforming formula of its own) -
# setting the DataBodyRange formula fills every row of the column and
# recalculates the dependent "Column1"/"Prev"/"Column2" columns too.
$col.DataBodyRange.Formula = "=ROUND([@Pos]*(5-5/5),0)"

# Restore the active selection to where the author left it.
$ws.Range("E28").Select()
